$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: header "time_taken", styled like the other header cells (copy E1's style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamps for data rows 2..31
$timestamps = @(
    "2021-10-05 10:52:06.835120",
    "2021-10-05 10:52:06.835135",
    "2021-10-05 10:52:06.835139",
    "2021-10-05 10:52:06.835142",
    "2021-10-05 10:52:06.835145",
    "2021-10-05 10:52:06.835148",
    "2021-10-05 10:52:06.835150",
    "2021-10-05 10:52:06.835153",
    "2021-10-05 10:52:06.835156",
    "2021-10-05 10:52:06.835159",
    "2021-10-05 10:52:06.835161",
    "2021-10-05 10:52:06.835164",
    "2021-10-05 10:52:06.835166",
    "2021-10-05 10:52:06.835169",
    "2021-10-05 10:52:06.835171",
    "2021-10-05 10:52:06.835174",
    "2021-10-05 10:52:06.835177",
    "2021-10-05 10:52:06.835180",
    "2021-10-05 10:52:06.835182",
    "2021-10-05 10:52:06.835185",
    "2021-10-05 10:52:06.835187",
    "2021-10-05 10:52:06.835190",
    "2021-10-05 10:52:06.835193",
    "2021-10-05 10:52:06.835195",
    "2021-10-05 10:52:06.835198",
    "2021-10-05 10:52:06.835201",
    "2021-10-05 10:52:06.835203",
    "2021-10-05 10:52:06.835206",
    "2021-10-05 10:52:06.835209",
    "2021-10-05 10:52:06.835211"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
